$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value, derived from the authoritative OOXML diff.
$changes = [ordered]@{
    "D2" = "288.54"
    "E2" = "1.18%"
    "D3" = "29.32"
    "E3" = "2.93%"
    "D4" = "5.099"
    "E4" = "2.63%"
    "D5" = "0.06686"
    "E5" = "3.10%"
    "D6" = "7.321"
    "E6" = "1.35%"
    "D7" = "3.408"
    "E7" = "1.28%"
    "D8" = "1.367"
    "E8" = "2.26%"
    "D9" = "0.9172"
    "E9" = "0.83%"
    "D10" = "0.1588"
    "E10" = "3.16%"
    "D11" = "0.06754"
    "E11" = "7.60%"
    "D12" = "0.07615"
    "E12" = "-0.72%"
    "D13" = "0.02940"
    "E13" = "-0.95%"
    "D14" = "0.08987"
    "E14" = "0.46%"
    "D15" = "0.001575"
    "E15" = "-1.55%"
    "D16" = "0.04504"
    "E16" = "0.73%"
    "D17" = "0.0006483"
    "E17" = "-0.83%"
    "D18" = "0.006291"
    "E18" = "3.92%"
    "D19" = "3.443"
    "E19" = "-0.56%"
    "D20" = "2.220"
    "E20" = "-0.95%"
    "E21" = "2.02%"
    "E22" = "-2.38%"
    "D23" = "4.063"
    "E23" = "1.58%"
    "D24" = "0.1582"
    "E24" = "1.89%"
    "D25" = "0.001192"
    "E25" = "0.51%"
    "D26" = "0.004114"
    "E26" = "-4.95%"
    "D27" = "0.0001200"
    "E27" = "1.61%"
    "D28" = "0.0001618"
    "E28" = "-1.04%"
    "D40" = "0.04243"
    "E40" = "1.98%"
    "D41" = "0.006727"
    "E41" = "0.60%"
    "D42" = "0.1239"
    "E42" = "0.58%"
    "D43" = "0.002249"
    "E43" = "5.06%"
    "D44" = "0.01342"
    "E44" = "13.96%"
    "D45" = "0.00005718"
    "E45" = "6.21%"
    "E46" = "1.81%"
    "D47" = "0.01307"
    "E47" = "-29.41%"
}

foreach ($addr in $changes.Keys) {
    $cell = $ws.Range($addr)
    # Force text interpretation so values like "288.54" or "1.18%" are
    # stored as literal strings (matching the source inlineStr cells)
    # instead of being auto-coerced into numbers/percentages.
    $cell.NumberFormat = "@"
    $cell.Value = $changes[$addr]
    # Drop back to the default style so no stray per-cell formatting sticks.
    $cell.Style = "Normal"
}

Write-Host "Applied $($changes.Count) cell updates"
